$wb = $excel.ActiveWorkbook

# Updated "want to go" counts (column F) on the "展览" sheet
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value = 301
$wsExhibition.Range("F7").Value = 1040
$wsExhibition.Range("F12").Value = 13282
$wsExhibition.Range("F13").Value = 163
$wsExhibition.Range("F14").Value = 10
$wsExhibition.Range("F16").Value = 5469
$wsExhibition.Range("F17").Value = 5567

# Same underlying events are duplicated on the "全部类型" aggregate sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F21").Value = 301
$wsAll.Range("F29").Value = 1040
$wsAll.Range("F34").Value = 13282
$wsAll.Range("F35").Value = 163
$wsAll.Range("F36").Value = 10
$wsAll.Range("F39").Value = 5469
$wsAll.Range("F40").Value = 5567
